# Applies the edit described by the diff:
#  - sets several new cell values (5) across rows 6, 15, 23, 24, 28
#  - updates the sheet view: frozen-pane scroll position and active selection

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6 ---
$ws.Range("M6").Value = 5
$ws.Range("N6").Value = 5
$ws.Range("O6").Value = 5
$ws.Range("V6").Value = 5
$ws.Range("W6").Value = 5
$ws.Range("Y6").Value = 5

# --- Row 15 ---
# M15 / N15 pick up the shaded "s=8" formatting used by the neighbouring
# I15:K15 cells, so copy that formatting across before writing the value.
$ws.Range("I15").Copy()
$ws.Range("M15:N15").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("H15").Value = 5
$ws.Range("M15").Value = 5
$ws.Range("N15").Value = 5
$ws.Range("O15").Value = 5
$ws.Range("R15").Value = 5
$ws.Range("U15").Value = 5
$ws.Range("V15").Value = 5
$ws.Range("X15").Value = 5
$ws.Range("Y15").Value = 5
$ws.Range("Z15").Value = 5

# --- Row 23 ---
$ws.Range("I23").Copy()
$ws.Range("M23:N23").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("M23").Value = 5
$ws.Range("N23").Value = 5
$ws.Range("O23").Value = 5
$ws.Range("Z23").Value = 5

# --- Row 24 ---
$ws.Range("Y24").Value = 5

# --- Row 28 ---
# V28 picks up the same shaded "s=16" formatting already used by U28.
$ws.Range("U24").Copy()
$ws.Range("V28").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("U28").Value = 5
$ws.Range("V28").Value = 5
$ws.Range("X28").Value = 5

# --- View / selection state ---
# Scroll the frozen pane so column H is the first visible column, then
# leave the cursor on H15 (mirrors the author's last selection/scroll).
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 8
$excel.ActiveWindow.ScrollRow = 3
$ws.Range("H15").Select()

Write-Host "Edit applied"
